$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 88.333336  # H12: 107 -> 88.333336
$ws.Cells.Item(12, 9).Value = 88.333336  # I12: 107 -> 88.333336
$ws.Cells.Item(12, 11).Value = 88.333336  # K12: 107 -> 88.333336
$ws.Cells.Item(12, 13).Value = 81.666664  # M12: 63 -> 81.666664
$ws.Cells.Item(49, 8).Value = 0  # H49: 3500 -> 0
$ws.Cells.Item(49, 10).Value = 0  # J49: 3500 -> 0
$ws.Cells.Item(49, 12).Value = 0  # L49: 10500 -> 0
$ws.Cells.Item(49, 14).Value = ""  # N49: delete (was DELETE)
$ws.Cells.Item(116, 8).Value = 28294.777  # H116: 30956.625 -> 28294.777
$ws.Cells.Item(116, 10).Value = 26201  # J116: 31001.25 -> 26201
$ws.Cells.Item(116, 12).Value = 26201  # L116: 31001.25 -> 26201
$ws.Cells.Item(116, 14).Value = -33085  # N116: -37885.25 -> -33085
$ws.Cells.Item(135, 8).Value = 2211.5  # H135: 2287.6316 -> 2211.5
$ws.Cells.Item(135, 9).Value = 1955.0588  # I135: 2029.4375 -> 1955.0588
$ws.Cells.Item(135, 11).Value = 17595.5292  # K135: 18264.9375 -> 17595.5292
$ws.Cells.Item(135, 13).Value = -15060.5292  # M135: -15729.9375 -> -15060.5292
$ws.Cells.Item(138, 8).Value = 10272.875  # H138: 10216 -> 10272.875
$ws.Cells.Item(138, 10).Value = 10505.15  # J138: 10431.798 -> 10505.15
$ws.Cells.Item(138, 12).Value = 31515.45  # L138: 31295.394 -> 31515.45
$ws.Cells.Item(138, 14).Value = -41795.45  # N138: -41575.394 -> -41795.45
$ws.Cells.Item(141, 8).Value = 5812.227  # H141: 6078.95 -> 5812.227
$ws.Cells.Item(141, 9).Value = 6015.95  # I141: 6334.9443 -> 6015.95
$ws.Cells.Item(141, 11).Value = 18047.85  # K141: 19004.8329 -> 18047.85
$ws.Cells.Item(141, 13).Value = -12867.85  # M141: -13824.8329 -> -12867.85

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 14012.226  # H32: 35000 -> 14012.226
$ws.Cells.Item(32, 9).Value = 13219.966  # I32: 0 -> 13219.966
$ws.Cells.Item(32, 10).Value = 25500  # J32: 35000 -> 25500
$ws.Cells.Item(32, 11).Value = 13219.966  # K32: 0 -> 13219.966
$ws.Cells.Item(32, 12).Value = 25500  # L32: 35000 -> 25500
$ws.Cells.Item(32, 13).Value = -12932.966  # M32: NEW -> -12932.966
$ws.Cells.Item(32, 14).Value = -26074  # N32: -35574 -> -26074
$ws.Cells.Item(37, 8).Value = 16749.75  # H37: 30500 -> 16749.75
$ws.Cells.Item(37, 10).Value = 57999  # J37: 58000 -> 57999
$ws.Cells.Item(37, 12).Value = 57999  # L37: 58000 -> 57999
$ws.Cells.Item(37, 14).Value = -58545  # N37: -58546 -> -58545
$ws.Cells.Item(61, 8).Value = 8060.143  # H61: 7971.357 -> 8060.143
$ws.Cells.Item(61, 9).Value = 6015.5293  # I61: 5859.95 -> 6015.5293
$ws.Cells.Item(61, 10).Value = 16749.75  # J61: 13249.875 -> 16749.75
$ws.Cells.Item(61, 11).Value = 6015.5293  # K61: 5859.95 -> 6015.5293
$ws.Cells.Item(61, 12).Value = 16749.75  # L61: 13249.875 -> 16749.75
$ws.Cells.Item(61, 13).Value = -5803.5293  # M61: -5647.95 -> -5803.5293
$ws.Cells.Item(61, 14).Value = -17173.75  # N61: -13673.875 -> -17173.75
$ws.Cells.Item(80, 8).Value = 84100  # H80: 86000 -> 84100
$ws.Cells.Item(80, 10).Value = 84100  # J80: 86000 -> 84100
$ws.Cells.Item(80, 12).Value = 84100  # L80: 86000 -> 84100
$ws.Cells.Item(80, 14).Value = -86096  # N80: -87996 -> -86096
$ws.Cells.Item(83, 8).Value = 84100  # H83: 86000 -> 84100
$ws.Cells.Item(83, 10).Value = 84100  # J83: 86000 -> 84100
$ws.Cells.Item(83, 12).Value = 252300  # L83: 258000 -> 252300
$ws.Cells.Item(83, 14).Value = -262284  # N83: -267984 -> -262284
$ws.Cells.Item(132, 8).Value = 4195.595  # H132: 4232.9517 -> 4195.595
$ws.Cells.Item(132, 9).Value = 3422.768  # I132: 3457 -> 3422.768
$ws.Cells.Item(132, 11).Value = 10268.304  # K132: 10371 -> 10268.304
$ws.Cells.Item(132, 13).Value = -7738.304  # M132: -7841 -> -7738.304
$ws.Cells.Item(136, 8).Value = 8060.143  # H136: 7971.357 -> 8060.143
$ws.Cells.Item(136, 9).Value = 6015.5293  # I136: 5859.95 -> 6015.5293
$ws.Cells.Item(136, 10).Value = 16749.75  # J136: 13249.875 -> 16749.75
$ws.Cells.Item(136, 11).Value = 18046.5879  # K136: 17579.85 -> 18046.5879
$ws.Cells.Item(136, 12).Value = 50249.25  # L136: 39749.625 -> 50249.25
$ws.Cells.Item(136, 13).Value = -15496.5879  # M136: -15029.85 -> -15496.5879
$ws.Cells.Item(136, 14).Value = -55349.25  # N136: -44849.625 -> -55349.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(35, 8).Value = 95520  # H35: 93828.57000000001 -> 95520
$ws.Cells.Item(35, 10).Value = 95520  # J35: 93828.57000000001 -> 95520
$ws.Cells.Item(35, 12).Value = 95520  # L35: 93828.57000000001 -> 95520
$ws.Cells.Item(35, 14).Value = -96140  # N35: -94448.57000000001 -> -96140
$ws.Cells.Item(86, 8).Value = 1310292.8  # H86: 1548182.2 -> 1310292.8
$ws.Cells.Item(86, 9).Value = 2127226.5  # I86: 2430858.8 -> 2127226.5
$ws.Cells.Item(86, 10).Value = 3198.8  # J86: 3498.5 -> 3198.8
$ws.Cells.Item(86, 11).Value = 2127226.5  # K86: 2430858.8 -> 2127226.5
$ws.Cells.Item(86, 12).Value = 3198.8  # L86: 3498.5 -> 3198.8
$ws.Cells.Item(86, 13).Value = -2126103.5  # M86: -2429735.8 -> -2126103.5
$ws.Cells.Item(86, 14).Value = -5444.8  # N86: -5744.5 -> -5444.8
$ws.Cells.Item(89, 8).Value = 1310292.8  # H89: 1548182.2 -> 1310292.8
$ws.Cells.Item(89, 9).Value = 2127226.5  # I89: 2430858.8 -> 2127226.5
$ws.Cells.Item(89, 10).Value = 3198.8  # J89: 3498.5 -> 3198.8
$ws.Cells.Item(89, 11).Value = 10636132.5  # K89: 12154294 -> 10636132.5
$ws.Cells.Item(89, 12).Value = 15994  # L89: 17492.5 -> 15994
$ws.Cells.Item(89, 13).Value = -10630516.5  # M89: -12148678 -> -10630516.5
$ws.Cells.Item(89, 14).Value = -27226  # N89: -28724.5 -> -27226
$ws.Cells.Item(92, 8).Value = 74000  # H92: 56000 -> 74000
$ws.Cells.Item(92, 10).Value = 74000  # J92: 56000 -> 74000
$ws.Cells.Item(92, 12).Value = 74000  # L92: 56000 -> 74000
$ws.Cells.Item(92, 14).Value = -78992  # N92: -60992 -> -78992

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 78799.92999999999  # H31: 84561.53999999999 -> 78799.92999999999
$ws.Cells.Item(31, 9).Value = 4433.25  # I31: 4481.8184 -> 4433.25
$ws.Cells.Item(31, 11).Value = 4433.25  # K31: 4481.8184 -> 4433.25
$ws.Cells.Item(31, 13).Value = -4138.25  # M31: -4186.8184 -> -4138.25
$ws.Cells.Item(34, 8).Value = 78799.92999999999  # H34: 84561.53999999999 -> 78799.92999999999
$ws.Cells.Item(34, 9).Value = 4433.25  # I34: 4481.8184 -> 4433.25
$ws.Cells.Item(34, 11).Value = 4433.25  # K34: 4481.8184 -> 4433.25
$ws.Cells.Item(34, 13).Value = -4231.25  # M34: -4279.8184 -> -4231.25
$ws.Cells.Item(50, 8).Value = 22995.055  # H50: 22939.5 -> 22995.055
$ws.Cells.Item(50, 10).Value = 56500  # J50: 56000 -> 56500
$ws.Cells.Item(50, 12).Value = 56500  # L50: 56000 -> 56500
$ws.Cells.Item(50, 14).Value = -57750  # N50: -57250 -> -57750
$ws.Cells.Item(58, 8).Value = 3049.64  # H58: 2997.7778 -> 3049.64
$ws.Cells.Item(58, 9).Value = 2926.3076  # I58: 2849.4 -> 2926.3076
$ws.Cells.Item(58, 11).Value = 2926.3076  # K58: 2849.4 -> 2926.3076
$ws.Cells.Item(58, 13).Value = -2723.3076  # M58: -2646.4 -> -2723.3076
$ws.Cells.Item(59, 8).Value = 40205  # H59: 42223 -> 40205
$ws.Cells.Item(59, 10).Value = 40205  # J59: 42223 -> 40205
$ws.Cells.Item(59, 12).Value = 40205  # L59: 42223 -> 40205
$ws.Cells.Item(59, 14).Value = -42495  # N59: -44513 -> -42495
$ws.Cells.Item(68, 8).Value = 68750  # H68: 0 -> 68750
$ws.Cells.Item(68, 10).Value = 68750  # J68: 0 -> 68750
$ws.Cells.Item(68, 12).Value = 68750  # L68: 0 -> 68750
$ws.Cells.Item(68, 14).Value = -70248  # N68: NEW -> -70248
$ws.Cells.Item(71, 8).Value = 68750  # H71: 0 -> 68750
$ws.Cells.Item(71, 10).Value = 68750  # J71: 0 -> 68750
$ws.Cells.Item(71, 12).Value = 206250  # L71: 0 -> 206250
$ws.Cells.Item(71, 14).Value = -213738  # N71: NEW -> -213738
$ws.Cells.Item(74, 8).Value = 99750  # H74: 92500 -> 99750
$ws.Cells.Item(74, 10).Value = 99750  # J74: 92500 -> 99750
$ws.Cells.Item(74, 12).Value = 99750  # L74: 92500 -> 99750
$ws.Cells.Item(74, 14).Value = -101498  # N74: -94248 -> -101498
$ws.Cells.Item(77, 8).Value = 99750  # H77: 92500 -> 99750
$ws.Cells.Item(77, 10).Value = 99750  # J77: 92500 -> 99750
$ws.Cells.Item(77, 12).Value = 299250  # L77: 277500 -> 299250
$ws.Cells.Item(77, 14).Value = -307986  # N77: -286236 -> -307986
$ws.Cells.Item(136, 8).Value = 3049.64  # H136: 2997.7778 -> 3049.64
$ws.Cells.Item(136, 9).Value = 2926.3076  # I136: 2849.4 -> 2926.3076
$ws.Cells.Item(136, 11).Value = 8778.9228  # K136: 8548.200000000001 -> 8778.9228
$ws.Cells.Item(136, 13).Value = -6228.9228  # M136: -5998.200000000001 -> -6228.9228

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 190.23077  # H2: 183.92592 -> 190.23077
$ws.Cells.Item(2, 9).Value = 72.09090999999999  # I2: 67.75 -> 72.09090999999999
$ws.Cells.Item(2, 11).Value = 432.5454599999999  # K2: 406.5 -> 432.5454599999999
$ws.Cells.Item(2, 13).Value = -319.5454599999999  # M2: -293.5 -> -319.5454599999999
$ws.Cells.Item(121, 8).Value = 15608053  # H121: 16351267 -> 15608053
$ws.Cells.Item(121, 9).Value = 1025  # I121: 1500 -> 1025
$ws.Cells.Item(121, 11).Value = 3075  # K121: 4500 -> 3075
$ws.Cells.Item(121, 13).Value = -1765  # M121: -3190 -> -1765

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(3, 8).Value = 6633.909  # H3: 6635.909 -> 6633.909
$ws.Cells.Item(3, 10).Value = 5638.7144  # J3: 5641.857 -> 5638.7144
$ws.Cells.Item(3, 12).Value = 5638.7144  # L3: 5641.857 -> 5638.7144
$ws.Cells.Item(3, 14).Value = -5870.7144  # N3: -5873.857 -> -5870.7144
$ws.Cells.Item(43, 8).Value = 20728.53  # H43: 23081.766 -> 20728.53
$ws.Cells.Item(43, 9).Value = 9748.75  # I43: 10090 -> 9748.75
$ws.Cells.Item(43, 10).Value = 47080  # J43: 46900 -> 47080
$ws.Cells.Item(43, 11).Value = 9748.75  # K43: 10090 -> 9748.75
$ws.Cells.Item(43, 12).Value = 47080  # L43: 46900 -> 47080
$ws.Cells.Item(43, 13).Value = -9597.75  # M43: -9939 -> -9597.75
$ws.Cells.Item(43, 14).Value = -47382  # N43: -47202 -> -47382
$ws.Cells.Item(46, 8).Value = 37090.547  # H46: 33166.332 -> 37090.547
$ws.Cells.Item(46, 9).Value = 0  # I46: 5000 -> 0
$ws.Cells.Item(46, 10).Value = 37090.547  # J46: 47249.5 -> 37090.547
$ws.Cells.Item(46, 11).Value = 0  # K46: 5000 -> 0
$ws.Cells.Item(46, 12).Value = 37090.547  # L46: 47249.5 -> 37090.547
$ws.Cells.Item(46, 13).Value = ""  # M46: delete (was DELETE)
$ws.Cells.Item(46, 14).Value = -37402.547  # N46: -47561.5 -> -37402.547
$ws.Cells.Item(132, 8).Value = 69071.3  # H132: 70529.39 -> 69071.3
$ws.Cells.Item(132, 9).Value = 5296.95  # I132: 5381.5127 -> 5296.95
$ws.Cells.Item(132, 11).Value = 15890.85  # K132: 16144.5381 -> 15890.85
$ws.Cells.Item(132, 13).Value = -13360.85  # M132: -13614.5381 -> -13360.85

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2675.0715  # H46: 2773.1538 -> 2675.0715
$ws.Cells.Item(46, 9).Value = 2313.7273  # I46: 2405.1 -> 2313.7273
$ws.Cells.Item(46, 11).Value = 2313.7273  # K46: 2405.1 -> 2313.7273
$ws.Cells.Item(46, 13).Value = -2125.7273  # M46: -2217.1 -> -2125.7273
$ws.Cells.Item(68, 8).Value = 3500  # H68: 0 -> 3500
$ws.Cells.Item(68, 9).Value = 3500  # I68: 0 -> 3500
$ws.Cells.Item(68, 11).Value = 3500  # K68: 0 -> 3500
$ws.Cells.Item(68, 13).Value = -2751  # M68: NEW -> -2751
$ws.Cells.Item(71, 8).Value = 3500  # H71: 0 -> 3500
$ws.Cells.Item(71, 9).Value = 3500  # I71: 0 -> 3500
$ws.Cells.Item(71, 11).Value = 17500  # K71: 0 -> 17500
$ws.Cells.Item(71, 13).Value = -13756  # M71: NEW -> -13756
